$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of trade data at row 9
$ws.Cells.Item(9, 1).Value = 8034.32
$ws.Cells.Item(9, 2).Value = 8256.42
$ws.Cells.Item(9, 3).Value = 19.36
$ws.Cells.Item(9, 4).Value = 18.84
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(9, 6).Value = -2.69

# Column G holds a date/time value, formatted like the other rows (style index 1 / numFmtId 22)
$ws.Cells.Item(9, 7).Value = 42612.673020833332
$ws.Cells.Item(9, 7).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(9, 8).Value = $false
